# Defect-Report.xlsx edit script
# Changes applied:
#  1. Insert a line break in the "15.Detailed Description" cell (D2) right
#     before "Version:128", while preserving the existing rich-text runs
#     (the underlined "Steps to reproduce the defect" heading and the
#     plain "step 1 - ..." instructions that follow it).
#  2. Update the active sheet selection to D2:D10 (the merged description
#     cell), matching the view state saved with the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$fullText = $cell.Text

$marker = "Version:128"
$markerIdx = $fullText.IndexOf($marker)

if ($markerIdx -ge 0) {
    # Characters() is 1-based; insert a newline immediately before the marker
    # using a zero-length character range at that position.
    $insertPos = $markerIdx + 1
    $insertChars = $cell.Characters($insertPos, 0)
    $insertChars.Text = "`n"

    # Re-resolve the text (it grew by one character) and restore the
    # original rich-text formatting for the two trailing runs that the
    # in-place edit may have flattened.
    $updatedText = $cell.Text

    $stepsHeading = "Steps to reproduce the defect                                           "
    $stepsIdx = $updatedText.IndexOf($stepsHeading)
    if ($stepsIdx -ge 0) {
        $headingChars = $cell.Characters($stepsIdx + 1, $stepsHeading.Length)
        $headingChars.Font.Name = "Calibri"
        $headingChars.Font.Size = 11
        $headingChars.Font.Underline = $true
    }

    $stepsBody = "step 1 - open browser and enter test url"
    $bodyIdx = $updatedText.IndexOf($stepsBody)
    if ($bodyIdx -ge 0) {
        $bodyLen = $updatedText.Length - $bodyIdx
        $bodyChars = $cell.Characters($bodyIdx + 1, $bodyLen)
        $bodyChars.Font.Name = "Calibri"
        $bodyChars.Font.Size = 11
        $bodyChars.Font.Underline = $false
    }
}

# Inserting the line break can make the engine mark row 2 with an explicit
# custom height; auto-fit it back so the row keeps behaving like the other
# default-height rows.
$ws.Rows.Item(2).AutoFit() | Out-Null

# Update the selection/active cell shown when the workbook is opened.
$ws.Range("D2:D10").Select() | Out-Null
